$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 40 (the "فرشه اسنان SENSODYNE" line item). This shifts all
#    subsequent product rows up by one, carrying their C/H/N/P (and other)
#    content with them.
$ws.Rows("40:40").Delete()

# 2. The "م" (sequence number) column is positional (row - 6) and must stay
#    that way after the shift - restore it for the affected rows.
for ($r = 40; $r -le 45; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 6
}

# 3. The former last line item has now moved out, and the new last row (45)
#    needs to become the new product: "معجون اسنان سيجنال 120 مل عرض".
$ws.Range("C45").Value2 = "معجون اسنان سيجنال 120 مل عرض"

# N45 / P45 hold their numbers as text. N45's style is already text-formatted,
# but P45's style is numeric ("0.00"), so Value2 would silently coerce the
# string to a Double unless we briefly swap it to a text format and back -
# that preserves the original style index (s="11") in the saved file.
$ws.Range("N45").Value2 = "70.00"

$origFmt = $ws.Range("P45").NumberFormat
$ws.Range("P45").NumberFormat = "@"
$ws.Range("P45").Value2 = "70.0000"
$ws.Range("P45").NumberFormat = $origFmt

# 4. Update the grand-total cell (now at P46 after the row shift).
$ws.Range("P46").Value2 = 1723.165

# 5. Update the generated timestamp in the footer (now row 47).
$ws.Range("A47").Value2 = "Tuesday, 23 September, 2025 5:05 PM"
